$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (C) and P_Value (D) columns with corrected values
$ws.Range("C2").Value = -1.278145547333787
$ws.Range("D2").Value = 0.2145195300955973

$ws.Range("C3").Value = -0.6410169888098358
$ws.Range("D3").Value = 0.5281334070391939

$ws.Range("C4").Value = -0.2498526058763735
$ws.Range("D4").Value = 0.8050195828821458

$ws.Range("C5").Value = 1.245832188492492
$ws.Range("D5").Value = 0.2259308812846104

$ws.Range("C6").Value = 0.3451876140220917
$ws.Range("D6").Value = 0.7332318798899893

$ws.Range("C7").Value = 0.6709979266217099
$ws.Range("D7").Value = 0.5092060350955552

$ws.Range("C8").Value = 1.854382719726822
$ws.Range("D8").Value = 0.07713876190455649
$ws.Range("G8").Value = "No"

$ws.Range("C9").Value = 0.465325117435665
$ws.Range("D9").Value = 0.6462719492043916

$ws.Range("C10").Value = 1.41200071442922
$ws.Range("D10").Value = 0.1719387008003486
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = 1.111373287106509
$ws.Range("D11").Value = 0.2784063462283677

$wb.Save()
